$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 500 (shifts existing rows 500..630 down to 501..631)
$ws.Rows.Item(500).Insert()

# Populate the newly inserted row 500 with the new data point
$ws.Cells.Item(500, 1).Value  = 5
$ws.Cells.Item(500, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(500, 3).Value  = "Maule"
$ws.Cells.Item(500, 4).Value  = 45135
$ws.Cells.Item(500, 5).Value  = 7
$ws.Cells.Item(500, 6).Value  = 100114014
$ws.Cells.Item(500, 7).Value  = "Betarraga"
$ws.Cells.Item(500, 8).Value  = "Sin especificar"
$ws.Cells.Item(500, 9).Value  = "Primera"
$ws.Cells.Item(500, 10).Value = 5000
$ws.Cells.Item(500, 11).Value = 500
$ws.Cells.Item(500, 12).Value = 500
$ws.Cells.Item(500, 13).Value = 500
$ws.Cells.Item(500, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(500, 15).Value = "Región del Maule"
$ws.Cells.Item(500, 16).Value = 100
$ws.Cells.Item(500, 17).Value = 5
$ws.Cells.Item(500, 18).Value = "Hortaliza"
